$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "20:02"
$ws.Range("F2").Value = 0.6458351956601358
$ws.Range("G2").Value = 98.52574419339665
$ws.Range("H2").Value = 99.81316567335617
